$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $cell = $ws.Range($cellAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '55.946.60'
Set-TextValue 'E2' '  -0.08%  '
Set-TextValue 'D3' '2.389.90'
Set-TextValue 'E3' '  -4.07%  '
Set-TextValue 'E4' '  +0.13%  '
Set-TextValue 'D5' '478.88'
Set-TextValue 'E5' '  -1.27%  '
Set-TextValue 'D6' '147.23'
Set-TextValue 'E6' '  +1.82%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  +0.23%  '
Set-TextValue 'E8' '  -1.63%  '
Set-TextValue 'D9' '2.389.41'
Set-TextValue 'E9' '  -4.61%  '
Set-TextValue 'D10' '0.0977'
Set-TextValue 'E10' '  +1.05%  '
Set-TextValue 'E11' '  -3.50%  '
Set-TextValue 'E12' '  -1.72%  '
Set-TextValue 'E13' '  +1.36%  '
Set-TextValue 'D14' '2.805.45'
Set-TextValue 'E14' '  -3.71%  '
Set-TextValue 'D15' '56.315.28'
Set-TextValue 'E15' '  +0.70%  '
Set-TextValue 'E16' '  -2.69%  '
Set-TextValue 'E17' '  -1.85%  '
Set-TextValue 'D18' '2.388.45'
Set-TextValue 'E18' '  -4.43%  '
Set-TextValue 'E19' '  +1.31%  '
Set-TextValue 'D20' '314.84'
Set-TextValue 'E20' '  -1.07%  '
Set-TextValue 'E21' '  -4.72%  '
Set-TextValue 'E22' '  -0.01%  '
Set-TextValue 'E23' '  -1.75%  '
Set-TextValue 'D24' '56.80'
Set-TextValue 'E24' '  -2.55%  '
Set-TextValue 'E25' '  +0.18%  '
Set-TextValue 'E26' '  -3.49%  '
Set-TextValue 'E27' '  -4.07%  '
Set-TextValue 'D28' '2.497.73'
Set-TextValue 'E28' '  -4.17%  '
Set-TextValue 'D29' '7.27'
Set-TextValue 'E29' '  -2.64%  '
Set-TextValue 'E30' '  -1.10%  '
Set-TextValue 'E31' '  +0.08%  '
Set-TextValue 'D32' '148.55'
Set-TextValue 'E32' '  +0.13%  '
Set-TextValue 'D33' '17.94'
Set-TextValue 'E33' '  -1.88%  '
Set-TextValue 'E34' '  -0.40%  '
Set-TextValue 'E35' '  -3.92%  '
Set-TextValue 'E36' '  -2.47%  '
Set-TextValue 'D37' '3.58'
Set-TextValue 'E37' '  -2.37%  '
Set-TextValue 'D38' '0.843'
Set-TextValue 'E38' '  -1.74%  '
Set-TextValue 'D39' '33.43'
Set-TextValue 'D40' '0.998'
Set-TextValue 'E40' '  +0.34%  '
Set-TextValue 'E41' '  +1.52%  '
Set-TextValue 'E42' '  -1.43%  '
Set-TextValue 'E43' '  -4.05%  '
Set-TextValue 'B44' 'Stellar'
Set-TextValue 'C44' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D44' '0.0946'
Set-TextValue 'E44' '  +5.07%  '
Set-TextValue 'B45' 'Mantle'
Set-TextValue 'C45' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D45' '0.583'
Set-TextValue 'E45' '  -4.35%  '
Set-TextValue 'D46' '10.21'
Set-TextValue 'E46' '  +0.47%  '
Set-TextValue 'B47' 'RenderToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D47' '4.63'
Set-TextValue 'E47' '  -2.31%  '
Set-TextValue 'B48' 'Bittensor'
Set-TextValue 'C48' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D48' '253.94'
Set-TextValue 'E48' '  -2.19%  '
Set-TextValue 'E49' '  -0.97%  '
Set-TextValue 'D50' '17.03'
Set-TextValue 'E50' '  -3.29%  '
Set-TextValue 'D51' '1.770.58'
Set-TextValue 'E51' '  -7.64%  '
